$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $systemParts = @()
        $otherParts = @()
        foreach ($p in $trimmed) {
            if ($p.CompareTo("System") -eq 0) {
                $systemParts += $p
            } else {
                $otherParts += $p
            }
        }
        $newParts = $otherParts + $systemParts
        $newVal = [string]::Join(", ", $newParts)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
